# Update the ozone MB-SHAP feature-importance table (Sheet1!A2:C103).
# The underlying SHAP values were recomputed (new run), which both
# changes every importance value in column C and re-sorts the rows
# by that new value (column A holds the stable original feature
# index, column B the feature label). A few labels for the
# categorical "type_0" .. "type_3" features are also now wrapped
# in "$...$" to match the LaTeX-style math formatting used by all
# the other feature labels.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item('Sheet1')

$ws.Cells.Item(2, 1).Value = 98
$ws.Cells.Item(2, 2).Value = '$(\langle pq \vert pq \rangle)_{3}$'
$ws.Cells.Item(2, 3).Value = [double]"0.0008892552225074907"
$ws.Cells.Item(3, 1).Value = 91
$ws.Cells.Item(3, 2).Value = '$(\langle pq \vert qp \rangle)_{2}$'
$ws.Cells.Item(3, 3).Value = [double]"0.0006776076343631585"
$ws.Cells.Item(4, 1).Value = 0
$ws.Cells.Item(4, 2).Value = '(h$_{p}$)$_{0}$'
$ws.Cells.Item(4, 3).Value = [double]"0.0004218182197812181"
$ws.Cells.Item(5, 1).Value = 26
$ws.Cells.Item(5, 2).Value = '$type_3$'
$ws.Cells.Item(5, 3).Value = [double]"0.0003765456554795044"
$ws.Cells.Item(6, 1).Value = 93
$ws.Cells.Item(6, 2).Value = '$(\langle rs \vert sr \rangle)_{2}$'
$ws.Cells.Item(6, 3).Value = [double]"0.0003169586951895845"
$ws.Cells.Item(7, 1).Value = 12
$ws.Cells.Item(7, 2).Value = 'h$_{q}$'
$ws.Cells.Item(7, 3).Value = [double]"0.0002821837277648635"
$ws.Cells.Item(8, 1).Value = 100
$ws.Cells.Item(8, 2).Value = '$(\langle rs\vert rs \rangle)_{3}$'
$ws.Cells.Item(8, 3).Value = [double]"0.000263599993350192"
$ws.Cells.Item(9, 1).Value = 20
$ws.Cells.Item(9, 2).Value = '(h$_{rs}$)$_{2}$'
$ws.Cells.Item(9, 3).Value = [double]"0.0001779199314815862"
$ws.Cells.Item(10, 1).Value = 7
$ws.Cells.Item(10, 2).Value = '(h$_{pq}$)$_{3}$'
$ws.Cells.Item(10, 3).Value = [double]"0.0001618232092465177"
$ws.Cells.Item(11, 1).Value = 78
$ws.Cells.Item(11, 2).Value = '$(\langle pq \vert rs \rangle)_{1}$'
$ws.Cells.Item(11, 3).Value = [double]"0.0001172200347640319"
$ws.Cells.Item(12, 1).Value = 97
$ws.Cells.Item(12, 2).Value = '$(\langle rr \vert rr \rangle)_{3}$'
$ws.Cells.Item(12, 3).Value = [double]"0.0001161616065326009"
$ws.Cells.Item(13, 1).Value = 85
$ws.Cells.Item(13, 2).Value = '$(\langle rs \vert sr \rangle)_{1}$'
$ws.Cells.Item(13, 3).Value = [double]"0.0001061750901893721"
$ws.Cells.Item(14, 1).Value = 2
$ws.Cells.Item(14, 2).Value = '(h$_{p}$)$_{2}$'
$ws.Cells.Item(14, 3).Value = [double]"9.137179801422224e-05"
$ws.Cells.Item(15, 1).Value = 3
$ws.Cells.Item(15, 2).Value = '(h$_{p}$)$_{3}$'
$ws.Cells.Item(15, 3).Value = [double]"8.436900978575903e-05"
$ws.Cells.Item(16, 1).Value = 90
$ws.Cells.Item(16, 2).Value = '$(\langle pq \vert pq \rangle)_{2}$'
$ws.Cells.Item(16, 3).Value = [double]"7.760257865862746e-05"
$ws.Cells.Item(17, 1).Value = 15
$ws.Cells.Item(17, 2).Value = '(h$_{r}$)$_{1}$'
$ws.Cells.Item(17, 3).Value = [double]"6.665416177921328e-05"
$ws.Cells.Item(18, 1).Value = 67
$ws.Cells.Item(18, 2).Value = '$(\eta_{r})_{3}$'
$ws.Cells.Item(18, 3).Value = [double]"6.609285404120526e-05"
$ws.Cells.Item(19, 1).Value = 92
$ws.Cells.Item(19, 2).Value = '$(\langle rs\vert rs \rangle)_{2}$'
$ws.Cells.Item(19, 3).Value = [double]"6.224987514591092e-05"
$ws.Cells.Item(20, 1).Value = 32
$ws.Cells.Item(20, 2).Value = '$F_{q}^{\text{SCF}}$'
$ws.Cells.Item(20, 3).Value = [double]"6.1428157110965e-05"
$ws.Cells.Item(21, 1).Value = 21
$ws.Cells.Item(21, 2).Value = '(h$_{rs}$)$_{3}$'
$ws.Cells.Item(21, 3).Value = [double]"6.031932114596838e-05"
$ws.Cells.Item(22, 1).Value = 52
$ws.Cells.Item(22, 2).Value = '$(F_{p}^{\text{SCF}})_{2}$'
$ws.Cells.Item(22, 3).Value = [double]"4.931637474349794e-05"
$ws.Cells.Item(23, 1).Value = 4
$ws.Cells.Item(23, 2).Value = '(h$_{pq}$)$_{0}$'
$ws.Cells.Item(23, 3).Value = [double]"4.677905199246879e-05"
$ws.Cells.Item(24, 1).Value = 6
$ws.Cells.Item(24, 2).Value = '(h$_{pq}$)$_{2}$'
$ws.Cells.Item(24, 3).Value = [double]"4.476402805394986e-05"
$ws.Cells.Item(25, 1).Value = 11
$ws.Cells.Item(25, 2).Value = '(h$_{pr}$)$_{3}$'
$ws.Cells.Item(25, 3).Value = [double]"4.441662147270612e-05"
$ws.Cells.Item(26, 1).Value = 22
$ws.Cells.Item(26, 2).Value = 'h$_{s}$'
$ws.Cells.Item(26, 3).Value = [double]"4.392259925439911e-05"
$ws.Cells.Item(27, 1).Value = 82
$ws.Cells.Item(27, 2).Value = '$(\langle pq \vert pq \rangle)_{1}$'
$ws.Cells.Item(27, 3).Value = [double]"4.008315689715553e-05"
$ws.Cells.Item(28, 1).Value = 16
$ws.Cells.Item(28, 2).Value = '(h$_{r}$)$_{2}$'
$ws.Cells.Item(28, 3).Value = [double]"3.841096886392836e-05"
$ws.Cells.Item(29, 1).Value = 42
$ws.Cells.Item(29, 2).Value = '$F_{s}$'
$ws.Cells.Item(29, 3).Value = [double]"3.832646188662629e-05"
$ws.Cells.Item(30, 1).Value = 36
$ws.Cells.Item(30, 2).Value = '$(F_{r}^{\text{SCF}})_{0}$'
$ws.Cells.Item(30, 3).Value = [double]"3.819268052465107e-05"
$ws.Cells.Item(31, 1).Value = 84
$ws.Cells.Item(31, 2).Value = '$(\langle rs\vert rs \rangle)_{1}$'
$ws.Cells.Item(31, 3).Value = [double]"3.761096707621735e-05"
$ws.Cells.Item(32, 1).Value = 17
$ws.Cells.Item(32, 2).Value = '(h$_{r}$)$_{3}$'
$ws.Cells.Item(32, 3).Value = [double]"3.741159163336254e-05"
$ws.Cells.Item(33, 1).Value = 40
$ws.Cells.Item(33, 2).Value = '$F_{s}^{\text{SCF}}$'
$ws.Cells.Item(33, 3).Value = [double]"3.595038121173624e-05"
$ws.Cells.Item(34, 1).Value = 77
$ws.Cells.Item(34, 2).Value = '$(\langle rs \vert sr \rangle)_{0}$'
$ws.Cells.Item(34, 3).Value = [double]"3.156820229723709e-05"
$ws.Cells.Item(35, 1).Value = 76
$ws.Cells.Item(35, 2).Value = '$(\langle rs\vert rs \rangle)_{0}$'
$ws.Cells.Item(35, 3).Value = [double]"2.953522764720236e-05"
$ws.Cells.Item(36, 1).Value = 73
$ws.Cells.Item(36, 2).Value = '$\langle ss \vert ss \rangle$'
$ws.Cells.Item(36, 3).Value = [double]"2.77964945386009e-05"
$ws.Cells.Item(37, 1).Value = 101
$ws.Cells.Item(37, 2).Value = '$(\langle rs \vert sr \rangle)_{3}$'
$ws.Cells.Item(37, 3).Value = [double]"2.739399323035937e-05"
$ws.Cells.Item(38, 1).Value = 71
$ws.Cells.Item(38, 2).Value = '$\langle qq \vert qq \rangle$'
$ws.Cells.Item(38, 3).Value = [double]"2.680716682388106e-05"
$ws.Cells.Item(39, 1).Value = 74
$ws.Cells.Item(39, 2).Value = '$(\langle pq \vert pq \rangle)_{0}$'
$ws.Cells.Item(39, 3).Value = [double]"2.674649015766396e-05"
$ws.Cells.Item(40, 1).Value = 99
$ws.Cells.Item(40, 2).Value = '$(\langle pq \vert qp \rangle)_{3}$'
$ws.Cells.Item(40, 3).Value = [double]"2.285409473116895e-05"
$ws.Cells.Item(41, 1).Value = 83
$ws.Cells.Item(41, 2).Value = '$(\langle pq \vert qp \rangle)_{1}$'
$ws.Cells.Item(41, 3).Value = [double]"2.149879209193023e-05"
$ws.Cells.Item(42, 1).Value = 86
$ws.Cells.Item(42, 2).Value = '$(\langle pq \vert rs \rangle)_{2}$'
$ws.Cells.Item(42, 3).Value = [double]"2.002463973982095e-05"
$ws.Cells.Item(43, 1).Value = 1
$ws.Cells.Item(43, 2).Value = '(h$_{p}$)$_{1}$'
$ws.Cells.Item(43, 3).Value = [double]"1.866916359273385e-05"
$ws.Cells.Item(44, 1).Value = 39
$ws.Cells.Item(44, 2).Value = '$(\eta_{r})_{0}$'
$ws.Cells.Item(44, 3).Value = [double]"1.838759804174708e-05"
$ws.Cells.Item(45, 1).Value = 5
$ws.Cells.Item(45, 2).Value = '(h$_{pq}$)$_{1}$'
$ws.Cells.Item(45, 3).Value = [double]"1.823433243068136e-05"
$ws.Cells.Item(46, 1).Value = 72
$ws.Cells.Item(46, 2).Value = '$(\langle rr \vert rr \rangle)_{0}$'
$ws.Cells.Item(46, 3).Value = [double]"1.78590309526702e-05"
$ws.Cells.Item(47, 1).Value = 13
$ws.Cells.Item(47, 2).Value = 'h$_{qs}$'
$ws.Cells.Item(47, 3).Value = [double]"1.716108678482271e-05"
$ws.Cells.Item(48, 1).Value = 10
$ws.Cells.Item(48, 2).Value = '(h$_{pr}$)$_{2}$'
$ws.Cells.Item(48, 3).Value = [double]"1.678683957782912e-05"
$ws.Cells.Item(49, 1).Value = 94
$ws.Cells.Item(49, 2).Value = '$(\langle pq \vert rs \rangle)_{3}$'
$ws.Cells.Item(49, 3).Value = [double]"1.616334270595508e-05"
$ws.Cells.Item(50, 1).Value = 34
$ws.Cells.Item(50, 2).Value = '$F_{q}$'
$ws.Cells.Item(50, 3).Value = [double]"1.614763434417835e-05"
$ws.Cells.Item(51, 1).Value = 89
$ws.Cells.Item(51, 2).Value = '$(\langle rr \vert rr \rangle)_{2}$'
$ws.Cells.Item(51, 3).Value = [double]"1.570083041490085e-05"
$ws.Cells.Item(52, 1).Value = 68
$ws.Cells.Item(52, 2).Value = '$(\langle pq \vert rs \rangle)_{0}$'
$ws.Cells.Item(52, 3).Value = [double]"1.508925831411055e-05"
$ws.Cells.Item(53, 1).Value = 38
$ws.Cells.Item(53, 2).Value = '$(F_{r})_{0}$'
$ws.Cells.Item(53, 3).Value = [double]"1.424880483408442e-05"
$ws.Cells.Item(54, 1).Value = 43
$ws.Cells.Item(54, 2).Value = '$\eta_{s}$'
$ws.Cells.Item(54, 3).Value = [double]"1.290858371315934e-05"
$ws.Cells.Item(55, 1).Value = 29
$ws.Cells.Item(55, 2).Value = '$(\omega_{p})_{0}$'
$ws.Cells.Item(55, 3).Value = [double]"1.2708536102171e-05"
$ws.Cells.Item(56, 1).Value = 80
$ws.Cells.Item(56, 2).Value = '$(\langle pp \vert pp \rangle)_{1}$'
$ws.Cells.Item(56, 3).Value = [double]"9.205525756419943e-06"
$ws.Cells.Item(57, 1).Value = 75
$ws.Cells.Item(57, 2).Value = '$(\langle pq \vert qp \rangle)_{0}$'
$ws.Cells.Item(57, 3).Value = [double]"9.081527878458481e-06"
$ws.Cells.Item(58, 1).Value = 9
$ws.Cells.Item(58, 2).Value = '(h$_{pr}$)$_{1}$'
$ws.Cells.Item(58, 3).Value = [double]"8.051452999969116e-06"
$ws.Cells.Item(59, 1).Value = 50
$ws.Cells.Item(59, 2).Value = '$(F_{r})_{1}$'
$ws.Cells.Item(59, 3).Value = [double]"7.358086587503526e-06"
$ws.Cells.Item(60, 1).Value = 60
$ws.Cells.Item(60, 2).Value = '$(F_{p}^{\text{SCF}})_{3}$'
$ws.Cells.Item(60, 3).Value = [double]"7.088282423278895e-06"
$ws.Cells.Item(61, 1).Value = 70
$ws.Cells.Item(61, 2).Value = '$(\langle pp \vert pp \rangle)_{0}$'
$ws.Cells.Item(61, 3).Value = [double]"6.868971842904003e-06"
$ws.Cells.Item(62, 1).Value = 8
$ws.Cells.Item(62, 2).Value = '(h$_{pr}$)$_{0}$'
$ws.Cells.Item(62, 3).Value = [double]"6.489261489953045e-06"
$ws.Cells.Item(63, 1).Value = 14
$ws.Cells.Item(63, 2).Value = '(h$_{r}$)$_{0}$'
$ws.Cells.Item(63, 3).Value = [double]"6.318327670216003e-06"
$ws.Cells.Item(64, 1).Value = 44
$ws.Cells.Item(64, 2).Value = '$(F_{p}^{\text{SCF}})_{1}$'
$ws.Cells.Item(64, 3).Value = [double]"6.188044321321396e-06"
$ws.Cells.Item(65, 1).Value = 88
$ws.Cells.Item(65, 2).Value = '$(\langle pp \vert pp \rangle)_{2}$'
$ws.Cells.Item(65, 3).Value = [double]"5.772977922155292e-06"
$ws.Cells.Item(66, 1).Value = 47
$ws.Cells.Item(66, 2).Value = '$(\eta_{p})_{1}$'
$ws.Cells.Item(66, 3).Value = [double]"5.525014456901111e-06"
$ws.Cells.Item(67, 1).Value = 24
$ws.Cells.Item(67, 2).Value = '$type_1$'
$ws.Cells.Item(67, 3).Value = [double]"4.81477978637754e-06"
$ws.Cells.Item(68, 1).Value = 19
$ws.Cells.Item(68, 2).Value = '(h$_{rs}$)$_{1}$'
$ws.Cells.Item(68, 3).Value = [double]"4.534229846952481e-06"
$ws.Cells.Item(69, 1).Value = 46
$ws.Cells.Item(69, 2).Value = '$(F_{p})_{1}$'
$ws.Cells.Item(69, 3).Value = [double]"4.276711829751058e-06"
$ws.Cells.Item(70, 1).Value = 64
$ws.Cells.Item(70, 2).Value = '$(F_{r}^{\text{SCF}})_{3}$'
$ws.Cells.Item(70, 3).Value = [double]"4.224029367340138e-06"
$ws.Cells.Item(71, 1).Value = 25
$ws.Cells.Item(71, 2).Value = '$type_2$'
$ws.Cells.Item(71, 3).Value = [double]"4.20460347654053e-06"
$ws.Cells.Item(72, 1).Value = 96
$ws.Cells.Item(72, 2).Value = '$(\langle pp \vert pp \rangle)_{3}$'
$ws.Cells.Item(72, 3).Value = [double]"3.74999930200561e-06"
$ws.Cells.Item(73, 1).Value = 59
$ws.Cells.Item(73, 2).Value = '$(\eta_{r})_{2}$'
$ws.Cells.Item(73, 3).Value = [double]"3.638890135239479e-06"
$ws.Cells.Item(74, 1).Value = 81
$ws.Cells.Item(74, 2).Value = '$(\langle rr \vert rr \rangle)_{1}$'
$ws.Cells.Item(74, 3).Value = [double]"3.629073264711836e-06"
$ws.Cells.Item(75, 1).Value = 56
$ws.Cells.Item(75, 2).Value = '$(F_{r}^{\text{SCF}})_{2}$'
$ws.Cells.Item(75, 3).Value = [double]"3.510968442180456e-06"
$ws.Cells.Item(76, 1).Value = 28
$ws.Cells.Item(76, 2).Value = '$(F_{p}^{\text{SCF}})_{0}$'
$ws.Cells.Item(76, 3).Value = [double]"2.92596904327747e-06"
$ws.Cells.Item(77, 1).Value = 66
$ws.Cells.Item(77, 2).Value = '$(F_{r})_{3}$'
$ws.Cells.Item(77, 3).Value = [double]"2.460464617639222e-06"
$ws.Cells.Item(78, 1).Value = 48
$ws.Cells.Item(78, 2).Value = '$(F_{r}^{\text{SCF}})_{1}$'
$ws.Cells.Item(78, 3).Value = [double]"2.369091502793501e-06"
$ws.Cells.Item(79, 1).Value = 31
$ws.Cells.Item(79, 2).Value = '$(\eta_{p})_{0}$'
$ws.Cells.Item(79, 3).Value = [double]"2.102365614181298e-06"
$ws.Cells.Item(80, 1).Value = 57
$ws.Cells.Item(80, 2).Value = '$(\omega_{r})_{2}$'
$ws.Cells.Item(80, 3).Value = [double]"1.77474324129859e-06"
$ws.Cells.Item(81, 1).Value = 54
$ws.Cells.Item(81, 2).Value = '$(F_{p})_{2}$'
$ws.Cells.Item(81, 3).Value = [double]"1.534202945859339e-06"
$ws.Cells.Item(82, 1).Value = 55
$ws.Cells.Item(82, 2).Value = '$(\eta_{p})_{2}$'
$ws.Cells.Item(82, 3).Value = [double]"1.506347614586752e-06"
$ws.Cells.Item(83, 1).Value = 63
$ws.Cells.Item(83, 2).Value = '$(\eta_{p})_{3}$'
$ws.Cells.Item(83, 3).Value = [double]"1.370469353465881e-06"
$ws.Cells.Item(84, 1).Value = 51
$ws.Cells.Item(84, 2).Value = '$(\eta_{r})_{1}$'
$ws.Cells.Item(84, 3).Value = [double]"7.4131753828976e-07"
$ws.Cells.Item(85, 1).Value = 58
$ws.Cells.Item(85, 2).Value = '$(F_{r})_{2}$'
$ws.Cells.Item(85, 3).Value = [double]"2.886832032721551e-07"
$ws.Cells.Item(86, 1).Value = 30
$ws.Cells.Item(86, 2).Value = '$(F_{p})_{0}$'
$ws.Cells.Item(86, 3).Value = [double]"8.471288768267502e-08"
$ws.Cells.Item(87, 1).Value = 95
$ws.Cells.Item(87, 2).Value = '$(\langle pq \vert sr \rangle)_{3}$'
$ws.Cells.Item(87, 3).Value = [double]"3.024944627609659e-08"
$ws.Cells.Item(88, 1).Value = 62
$ws.Cells.Item(88, 2).Value = '$(F_{p})_{3}$'
$ws.Cells.Item(88, 3).Value = [double]"2.897218793969686e-08"
$ws.Cells.Item(89, 1).Value = 23
$ws.Cells.Item(89, 2).Value = '$type_0$'
$ws.Cells.Item(89, 3).Value = [double]"2.806609936770687e-08"
$ws.Cells.Item(90, 1).Value = 69
$ws.Cells.Item(90, 2).Value = '$(\langle pq \vert sr \rangle)_{0}$'
$ws.Cells.Item(90, 3).Value = [double]"2.718529393792282e-08"
$ws.Cells.Item(91, 1).Value = 37
$ws.Cells.Item(91, 2).Value = '$(\omega_{r})_{0}$'
$ws.Cells.Item(91, 3).Value = [double]"2.129332987139472e-08"
$ws.Cells.Item(92, 1).Value = 79
$ws.Cells.Item(92, 2).Value = '$(\langle pq \vert sr \rangle)_{1}$'
$ws.Cells.Item(92, 3).Value = [double]"2.08494733057455e-08"
$ws.Cells.Item(93, 1).Value = 45
$ws.Cells.Item(93, 2).Value = '$(\omega_{p})_{1}$'
$ws.Cells.Item(93, 3).Value = [double]"1.967055429854626e-08"
$ws.Cells.Item(94, 1).Value = 87
$ws.Cells.Item(94, 2).Value = '$(\langle pq \vert sr \rangle)_{2}$'
$ws.Cells.Item(94, 3).Value = [double]"1.925970954975286e-08"
$ws.Cells.Item(95, 1).Value = 65
$ws.Cells.Item(95, 2).Value = '$(\omega_{r})_{3}$'
$ws.Cells.Item(95, 3).Value = [double]"1.892729815990517e-08"
$ws.Cells.Item(96, 1).Value = 41
$ws.Cells.Item(96, 2).Value = '$\omega_{s}$'
$ws.Cells.Item(96, 3).Value = [double]"1.815645642307919e-08"
$ws.Cells.Item(97, 1).Value = 49
$ws.Cells.Item(97, 2).Value = '$(\omega_{r})_{1}$'
$ws.Cells.Item(97, 3).Value = [double]"1.484958757914115e-08"
$ws.Cells.Item(98, 1).Value = 53
$ws.Cells.Item(98, 2).Value = '$(\omega_{p})_{2}$'
$ws.Cells.Item(98, 3).Value = [double]"1.20669038860983e-08"
$ws.Cells.Item(99, 1).Value = 18
$ws.Cells.Item(99, 2).Value = '(h$_{rs}$)$_{0}$'
$ws.Cells.Item(99, 3).Value = [double]"1.089875441951313e-08"
$ws.Cells.Item(100, 1).Value = 35
$ws.Cells.Item(100, 2).Value = '$\eta_{q}$'
$ws.Cells.Item(100, 3).Value = [double]"3.400179300333807e-09"
$ws.Cells.Item(101, 1).Value = 27
$ws.Cells.Item(101, 2).Value = '$\mathbf{b}$'
$ws.Cells.Item(101, 3).Value = [double]"2.685787764750484e-09"
$ws.Cells.Item(102, 1).Value = 61
$ws.Cells.Item(102, 2).Value = '$(\omega_{p})_{3}$'
$ws.Cells.Item(102, 3).Value = [double]"1.61497801443314e-09"
$ws.Cells.Item(103, 1).Value = 33
$ws.Cells.Item(103, 2).Value = '$\omega_{q}$'
$ws.Cells.Item(103, 3).Value = [double]"0"
